$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, copying the header formatting used by the
# other header cells (e.g. G1) since named cell Styles aren't used here -
# the formatting comes from direct cell formatting (bold font, border,
# centered alignment).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" column values for the existing data rows.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
